# "more adjustments to ACH" -- append new check-register rows 377-382
# to the "Check Register" sheet (ACH / check payments for 6/25 and 6/30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Check Register")

# --- Rows 377-379 (check run dated 2024-06-25 / serial 45468) ---
# These pick up the same formatting already used by the column (left-aligned
# text for Check #/Payee/Cash Account, m/d/yy date, and a right-aligned
# money format), same as the rest of the register.

$ws.Cells.Item(377, 1).Value = "12929"
$ws.Cells.Item(377, 2).Value = 45468
$ws.Cells.Item(377, 3).Value = "Advance Scale"
$ws.Cells.Item(377, 4).Value = "11040"
$ws.Cells.Item(377, 5).Value = 3128.06

$ws.Cells.Item(378, 1).Value = "12930"
$ws.Cells.Item(378, 2).Value = 45468
$ws.Cells.Item(378, 3).Value = "Employment Screening Services, Inc"
$ws.Cells.Item(378, 4).Value = "11040"
$ws.Cells.Item(378, 5).Value = 300.98

$ws.Cells.Item(379, 1).Value = "12931"
$ws.Cells.Item(379, 2).Value = 45468
$ws.Cells.Item(379, 3).Value = "Linemark"
$ws.Cells.Item(379, 4).Value = "11040"
$ws.Cells.Item(379, 5).Value = 2404.5

# --- Rows 380-382 (check run dated 2024-06-30 / serial 45473) ---
# Entered in a later session, so Excel registers them against a
# freshly-derived (but visually identical) style: left aligned text /
# date / right aligned money, explicitly carrying the cell protection
# (locked) state along with the rest of the formatting.

$ws.Cells.Item(380, 1).Value = "12932"
$ws.Cells.Item(380, 2).Value = 45473
$ws.Cells.Item(380, 3).Value = "Linemark"
$ws.Cells.Item(380, 4).Value = "11040"
$ws.Cells.Item(380, 5).Value = 97.32

$ws.Cells.Item(381, 1).Value = "12933"
$ws.Cells.Item(381, 2).Value = 45473
$ws.Cells.Item(381, 3).Value = "Neptune and Co., Inc."
$ws.Cells.Item(381, 4).Value = "11040"
$ws.Cells.Item(381, 5).Value = 10200

$ws.Cells.Item(382, 1).Value = "12934"
$ws.Cells.Item(382, 2).Value = 45473
$ws.Cells.Item(382, 3).Value = "Office Equipment Svcs"
$ws.Cells.Item(382, 4).Value = "11040"
$ws.Cells.Item(382, 5).Value = 149.95

$newRows = $ws.Range("A380:E382")
$newRows.Locked = $true

# --- view state: scroll the frozen pane down and select the newly
# entered block, matching where the user left off editing ---
$ws.Activate()
$ws.Range("A380:E382").Select()
